# Regenerate merged AHB files
#
# 1. Rename the header row: "<name>_old" -> "<name>_FV2404" and
#    "<name>_new" -> "<name>_FV2410" (the sheet diffs the FV2404 and
#    FV2410 Anwendungshandbuch versions against each other).
# 2. Turn the used range A1:U62 into a native Excel Table ("Table1") with
#    its AutoFilter, so the merged sheet can be filtered/sorted per column.
# 3. Freeze the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- 1. header renames -----------------------------------------------
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = [string]$cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2404"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2410"
        }
    }
}

# --- 2. convert the range into a table --------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A100:U100")

# Stash the header row's existing formatting out of the way first: Excel
# captures any pre-existing header formatting that differs from the new
# table style into a dxf override, which we don't want here. Restoring
# the identical formatting afterwards keeps the header's look unchanged.
$headerRange.Copy()
$scratchRange.PasteSpecial(-4122) # xlPasteFormats
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U62")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(100).Delete()

# --- 3. freeze the header row ------------------------------------------
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
